$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-03-29 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-03-30 Saturday", 2)

# Update the multiplication problems in the table, cell by cell (row, column)
# so that duplicate expressions (e.g. "978x4=" appearing in two different
# cells) are each replaced with their own distinct target value.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "892×6="
$t.Cell(1, 2).Range.Text  = "132×2="
$t.Cell(1, 3).Range.Text  = "729×3="
$t.Cell(1, 4).Range.Text  = "126×7="
$t.Cell(1, 5).Range.Text  = "393×8="

$t.Cell(5, 1).Range.Text  = "660×7="
$t.Cell(5, 2).Range.Text  = "417×8="
$t.Cell(5, 3).Range.Text  = "271×4="
$t.Cell(5, 4).Range.Text  = "892×4="
$t.Cell(5, 5).Range.Text  = "444×4="

$t.Cell(10, 1).Range.Text = "363×2="
$t.Cell(10, 2).Range.Text = "326×4="
$t.Cell(10, 3).Range.Text = "185×5="
$t.Cell(10, 4).Range.Text = "370×2="
$t.Cell(10, 5).Range.Text = "613×7="

$t.Cell(15, 1).Range.Text = "461×9="
$t.Cell(15, 2).Range.Text = "964×5="
$t.Cell(15, 3).Range.Text = "123×7="
$t.Cell(15, 4).Range.Text = "493×5="
$t.Cell(15, 5).Range.Text = "284×7="

$t.Cell(20, 1).Range.Text = "389×6="
$t.Cell(20, 2).Range.Text = "146×3="
$t.Cell(20, 3).Range.Text = "463×6="
$t.Cell(20, 4).Range.Text = "183×2="
$t.Cell(20, 5).Range.Text = "517×7="
